# "add prices to items"
# Adds a new "Price" column (C) to the items list, with values 1-10 for the
# 10 existing rows, formatted as Naira currency, and nudges the sheet's
# selection/page setup to match the saved workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in C1
$ws.Range("C1").Value = "Price"

# Price values for the 10 item rows (rows 2-11): 1, 2, 3, ... 10
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $i + 1
}

# Apply Naira currency number format to the new Price column's data cells
$ws.Range("C2:C11").NumberFormat = "[$₦-470]#,##0.00"

# Match the saved selection state
[void]$ws.Range("J4").Select()

# Match the saved page orientation
$ws.PageSetup.Orientation = 1
